$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.922
$ws.Range("B3").Value = 6.798999999999999
$ws.Range("B5").Value = 6.243
$ws.Range("C7").Value = -12.871
$ws.Range("A9").Value = -21.22
$ws.Range("C9").Value = -12.449
$ws.Range("B11").Value = 6.568
$ws.Range("B12").Value = 6.404999999999999
$ws.Range("A13").Value = -21.786
$ws.Range("A16").Value = -20.763
$ws.Range("A18").Value = -21.751
$ws.Range("A20").Value = -21.664
$ws.Range("B21").Value = 6.568000000000001
$ws.Range("C21").Value = -13.339
